# Update forecast values on the "Forecast Comparison" sheet (auto-arima
# column removed from the model mix, so Amazon Mean/P70/P80/P90 forecasts
# were recalculated for weeks W01-W16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Row -> D (Amazon Mean), E (Amazon P70), F (Amazon P80), G (Amazon P90)
$values = @{
    2  = @(20, 25, 29, 37)
    3  = @(16, 20, 24, 31)
    4  = @(15, 19, 23, 29)
    5  = @(17, 20, 25, 32)
    6  = @(17, 21, 25, 33)
    7  = @(16, 20, 24, 32)
    8  = @(17, 21, 26, 34)
    9  = @(17, 21, 27, 35)
    10 = @(16, 20, 25, 33)
    11 = @(17, 20, 26, 34)
    12 = @(17, 21, 26, 36)
    13 = @(18, 22, 28, 38)
    14 = @(16, 20, 26, 35)
    15 = @(16, 19, 25, 35)
    16 = @(16, 19, 25, 35)
    17 = @(15, 19, 24, 34)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 4).Value = $rowValues[0]  # D - Amazon Mean Forecast
    $ws.Cells.Item($row, 5).Value = $rowValues[1]  # E - Amazon P70 Forecast
    $ws.Cells.Item($row, 6).Value = $rowValues[2]  # F - Amazon P80 Forecast
    $ws.Cells.Item($row, 7).Value = $rowValues[3]  # G - Amazon P90 Forecast
}

$wb.Save()
